$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: title line
#   " – Assignment 1 exercise"  ->  5 separate runs (same bold/sz28
#   formatting) spelling out " – Assignment 2 Exercise"
# ------------------------------------------------------------------

$titleFind = $d.Content.Find
$titleFind.Execute(" – Assignment 1 exercise") | Out-Null
if (-not $titleFind.Found) {
    throw "Could not locate the title run to replace"
}
$titleRange = $titleFind.Parent

$titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> – Assignment </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>2</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>E</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>xercise</w:t></w:r>' +
    '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$titleRange.InsertXML($titleXml)

# ------------------------------------------------------------------
# Change 2: fill in the trailing empty paragraph with the
#   "Version management" commentary paragraph (15 runs)
# ------------------------------------------------------------------

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$bodyXml = '<w:r><w:t xml:space="preserve">One important aspect of Version management is </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Version control </w:t></w:r>' +
    '<w:r><w:t>as</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> it allows</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> one project to be worked on by</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> multiple people simultaneously,</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> whilst</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> keep</w:t></w:r>' +
    '<w:r><w:t>ing</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> track of changes over time</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> allowing them to</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> collaborate more efficiently, and</w:t></w:r>' +
    '<w:r><w:t>, in the case of errors or mistakes, quickly</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> revert to a previous </w:t></w:r>' +
    '<w:r><w:t>version.</w:t></w:r>'

$paraXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $bodyXml + '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($paraXml)

# InsertXML on a range that covered the whole (empty) paragraph mark
# inserts the new paragraph *before* it and leaves the old, now
# redundant, empty paragraph behind - remove that leftover paragraph.
$countAfterInsert = $d.Paragraphs.Count
$newLastPara = $d.Paragraphs.Item($countAfterInsert)
$priorPara = $d.Paragraphs.Item($countAfterInsert - 1)
$cleanupRange = $d.Range($priorPara.Range.End - 1, $newLastPara.Range.End)
$cleanupRange.Delete()
